$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee Info")

$ws.Range("A12").Value = ""
$ws.Range("B12").Value = "2"
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = "Mezzanine"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
